$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, -0.04529952735016188, 0.3807916908414835),
    @(3, 0.1863946186624733, 0.04870315033453594),
    @(4, -0.1979653412675272, -0.04578366880077607),
    @(5, -0.04628039727168934, 0.0527279543444332),
    @(6, 0.03697795705066301, -0.2111855394777243),
    @(7, 0.2508200006377966, 0.1958649351207464),
    @(8, -0.03976996254254109, 0.5165749593323006),
    @(9, 0.5275972255920272, -0.2589894288161541),
    @(10, -0.09331759947224851, 0.2269473968572833),
    @(11, 0.3274540970952823, 0.2278432754087099),
    @(12, 0.2519163786991722, 0.005729994812295732),
    @(13, -0.09040269537094542, 0.3040341500528089),
    @(14, 0.3948949280184769, -0.01068316406921473),
    @(15, 0.3356224698231056, -0.1842273898718651),
    @(16, 0.1776075311436495, 0.2458715389224503),
    @(17, 0.299927611655373, 0.3747902625997055)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

$wb.Save()
